# corrección graficas y preguntas
# falata corregir las descripciones
#
# Update the demographic counts in row 2 and move the selection to N2
# (with the view scrolled so column G is the first visible column),
# matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 data corrections ---
$ws.Range("A2").Value = 25
$ws.Range("B2").Value = 28
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 44
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 26
$ws.Range("I2").Value = 23
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 3

# --- View / selection state ---
# Scroll so column G is left-most visible, then select N2 (matches the
# saved sheetView topLeftCell="G1" / selection activeCell="N2" sqref="N2").
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N2").Select()
